# Scheduled Sheets runner: refresh Market Board averages/profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the Leve
# crafting-class workbooks (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 356.625
$ws.Cells.Item(33, 10).Value = 750
$ws.Cells.Item(33, 12).Value = 750
$ws.Cells.Item(33, 14).Value = -1208

$ws.Cells.Item(42, 8).Value = 330.35715
$ws.Cells.Item(42, 9).Value = 317.81818
$ws.Cells.Item(42, 10).Value = 376.33334
$ws.Cells.Item(42, 11).Value = 953.45454
$ws.Cells.Item(42, 12).Value = 1129.00002
$ws.Cells.Item(42, 13).Value = -723.45454
$ws.Cells.Item(42, 14).Value = -1589.00002

$ws.Cells.Item(43, 8).Value = 3740.5557
$ws.Cells.Item(43, 9).Value = 2228.875
$ws.Cells.Item(43, 11).Value = 2228.875
$ws.Cells.Item(43, 13).Value = -2159.875

$ws.Cells.Item(76, 8).Value = 4450.7
$ws.Cells.Item(76, 10).Value = 5999
$ws.Cells.Item(76, 12).Value = 5999
$ws.Cells.Item(76, 14).Value = -6629

$ws.Cells.Item(79, 8).Value = 4450.7
$ws.Cells.Item(79, 10).Value = 5999
$ws.Cells.Item(79, 12).Value = 5999
$ws.Cells.Item(79, 14).Value = -8183

$ws.Cells.Item(86, 8).Value = 1794.3334
$ws.Cells.Item(86, 9).Value = 1537.2667
$ws.Cells.Item(86, 11).Value = 1537.2667
$ws.Cells.Item(86, 13).Value = -414.2666999999999

$ws.Cells.Item(89, 8).Value = 1794.3334
$ws.Cells.Item(89, 9).Value = 1537.2667
$ws.Cells.Item(89, 11).Value = 7686.3335
$ws.Cells.Item(89, 13).Value = -2070.3335

$ws.Cells.Item(111, 8).Value = 192.66667
$ws.Cells.Item(111, 9).Value = 192.66667
$ws.Cells.Item(111, 11).Value = 578.00001
$ws.Cells.Item(111, 13).Value = 2488.99999

$ws.Cells.Item(113, 8).Value = 4768.625
$ws.Cells.Item(113, 9).Value = 2135
$ws.Cells.Item(113, 10).Value = 6348.8
$ws.Cells.Item(113, 11).Value = 2135
$ws.Cells.Item(113, 12).Value = 6348.8
$ws.Cells.Item(113, 13).Value = 1119
$ws.Cells.Item(113, 14).Value = -12856.8

$ws.Cells.Item(131, 8).Value = 2462.739
$ws.Cells.Item(131, 9).Value = 2107.25
$ws.Cells.Item(131, 10).Value = 4832.6665
$ws.Cells.Item(131, 11).Value = 6321.75
$ws.Cells.Item(131, 12).Value = 14497.9995
$ws.Cells.Item(131, 13).Value = -1281.75
$ws.Cells.Item(131, 14).Value = -24577.9995

$ws.Cells.Item(132, 8).Value = 19928.309
$ws.Cells.Item(132, 9).Value = 22465.734
$ws.Cells.Item(132, 10).Value = 2673.8
$ws.Cells.Item(132, 11).Value = 67397.202
$ws.Cells.Item(132, 12).Value = 8021.400000000001
$ws.Cells.Item(132, 13).Value = -64867.202
$ws.Cells.Item(132, 14).Value = -13081.4

$ws.Cells.Item(135, 8).Value = 1186.0555
$ws.Cells.Item(135, 9).Value = 584.3125
$ws.Cells.Item(135, 11).Value = 5258.8125
$ws.Cells.Item(135, 13).Value = -2723.8125

$ws.Cells.Item(137, 8).Value = 14530.292
$ws.Cells.Item(137, 9).Value = 18296.223
$ws.Cells.Item(137, 11).Value = 54888.66900000001
$ws.Cells.Item(137, 13).Value = -52338.66900000001


$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 299.5
$ws.Cells.Item(4, 9).Value = 299.5
$ws.Cells.Item(4, 11).Value = 299.5
$ws.Cells.Item(4, 13).Value = -183.5

$ws.Cells.Item(28, 8).Value = 5913.4443
$ws.Cells.Item(28, 9).Value = 5913.4443
$ws.Cells.Item(28, 11).Value = 5913.4443
$ws.Cells.Item(28, 13).Value = -5721.4443

$ws.Cells.Item(45, 8).Value = 3751.4375
$ws.Cells.Item(45, 9).Value = 2442.25
$ws.Cells.Item(45, 11).Value = 2442.25
$ws.Cells.Item(45, 13).Value = -2065.25

$ws.Cells.Item(88, 8).Value = 7844
$ws.Cells.Item(88, 9).Value = 4465.6665
$ws.Cells.Item(88, 11).Value = 4465.6665
$ws.Cells.Item(88, 13).Value = -4059.6665

$ws.Cells.Item(91, 8).Value = 7844
$ws.Cells.Item(91, 9).Value = 4465.6665
$ws.Cells.Item(91, 11).Value = 4465.6665
$ws.Cells.Item(91, 13).Value = -3061.6665

$ws.Cells.Item(99, 8).Value = 5913.4443
$ws.Cells.Item(99, 9).Value = 5913.4443
$ws.Cells.Item(99, 11).Value = 5913.4443
$ws.Cells.Item(99, 13).Value = -2918.4443

$ws.Cells.Item(122, 8).Value = 1690.68
$ws.Cells.Item(122, 9).Value = 1583.9048
$ws.Cells.Item(122, 11).Value = 4751.7144
$ws.Cells.Item(122, 13).Value = -2301.7144


$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 938
$ws.Cells.Item(22, 9).Value = 906
$ws.Cells.Item(22, 10).Value = 1002
$ws.Cells.Item(22, 11).Value = 906
$ws.Cells.Item(22, 12).Value = 1002
$ws.Cells.Item(22, 13).Value = -733
$ws.Cells.Item(22, 14).Value = -1348

$ws.Cells.Item(36, 8).Value = 2924.75
$ws.Cells.Item(36, 9).Value = 2933
$ws.Cells.Item(36, 10).Value = 2900
$ws.Cells.Item(36, 11).Value = 2933
$ws.Cells.Item(36, 12).Value = 2900
$ws.Cells.Item(36, 13).Value = -2399
$ws.Cells.Item(36, 14).Value = -3968

$ws.Cells.Item(86, 8).Value = 4478.8335
$ws.Cells.Item(86, 9).Value = 1749.5
$ws.Cells.Item(86, 11).Value = 1749.5
$ws.Cells.Item(86, 13).Value = -626.5

$ws.Cells.Item(89, 8).Value = 4478.8335
$ws.Cells.Item(89, 9).Value = 1749.5
$ws.Cells.Item(89, 11).Value = 8747.5
$ws.Cells.Item(89, 13).Value = -3131.5

$ws.Cells.Item(94, 8).Value = 5105.1924
$ws.Cells.Item(94, 9).Value = 5612.4
$ws.Cells.Item(94, 10).Value = 3414.5
$ws.Cells.Item(94, 11).Value = 5612.4
$ws.Cells.Item(94, 12).Value = 3414.5
$ws.Cells.Item(94, 13).Value = -5161.4
$ws.Cells.Item(94, 14).Value = -4316.5

$ws.Cells.Item(99, 8).Value = 2142.2856
$ws.Cells.Item(99, 9).Value = 2071
$ws.Cells.Item(99, 10).Value = 2320.5
$ws.Cells.Item(99, 11).Value = 2071
$ws.Cells.Item(99, 12).Value = 2320.5
$ws.Cells.Item(99, 13).Value = -573
$ws.Cells.Item(99, 14).Value = -5316.5


$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1039.8
$ws.Cells.Item(22, 9).Value = 749
$ws.Cells.Item(22, 11).Value = 749
$ws.Cells.Item(22, 13).Value = -399

$ws.Cells.Item(122, 8).Value = 2502.5454
$ws.Cells.Item(122, 9).Value = 2679.875
$ws.Cells.Item(122, 10).Value = 2029.6666
$ws.Cells.Item(122, 11).Value = 8039.625
$ws.Cells.Item(122, 12).Value = 6088.9998
$ws.Cells.Item(122, 13).Value = -5589.625
$ws.Cells.Item(122, 14).Value = -10988.9998

$ws.Cells.Item(132, 8).Value = 40851.88
$ws.Cells.Item(132, 9).Value = 48476.57
$ws.Cells.Item(132, 10).Value = 822.25
$ws.Cells.Item(132, 11).Value = 145429.71
$ws.Cells.Item(132, 12).Value = 2466.75
$ws.Cells.Item(132, 13).Value = -142899.71
$ws.Cells.Item(132, 14).Value = -7526.75


$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(25, 8).Value = 1600.2
$ws.Cells.Item(25, 9).Value = 2001
$ws.Cells.Item(25, 10).Value = 1500
$ws.Cells.Item(25, 11).Value = 6003
$ws.Cells.Item(25, 12).Value = 4500
$ws.Cells.Item(25, 13).Value = -5834
$ws.Cells.Item(25, 14).Value = -4838

$ws.Cells.Item(30, 8).Value = 1600.2
$ws.Cells.Item(30, 9).Value = 2001
$ws.Cells.Item(30, 10).Value = 1500
$ws.Cells.Item(30, 11).Value = 6003
$ws.Cells.Item(30, 12).Value = 4500
$ws.Cells.Item(30, 13).Value = -5901
$ws.Cells.Item(30, 14).Value = -4704

$ws.Cells.Item(47, 8).Value = 3314.7144
$ws.Cells.Item(47, 9).Value = 3314.7144
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 11).Value = 9944.143199999999
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 13).Value = -9513.143199999999
$ws.Cells.Item(47, 14).ClearContents()

$ws.Cells.Item(122, 8).Value = 820.4666999999999
$ws.Cells.Item(122, 9).Value = 440
$ws.Cells.Item(122, 10).Value = 915.5833
$ws.Cells.Item(122, 11).Value = 3960
$ws.Cells.Item(122, 12).Value = 8240.2497
$ws.Cells.Item(122, 13).Value = -1510
$ws.Cells.Item(122, 14).Value = -13140.2497

$ws.Cells.Item(129, 8).Value = 2235.5
$ws.Cells.Item(129, 9).Value = 1707.1818
$ws.Cells.Item(129, 11).Value = 5121.5454
$ws.Cells.Item(129, 13).Value = -121.5454

$ws.Cells.Item(131, 8).Value = 429964.8
$ws.Cells.Item(131, 10).Value = 3249.25
$ws.Cells.Item(131, 12).Value = 9747.75
$ws.Cells.Item(131, 14).Value = -19827.75


$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 8098.3
$ws.Cells.Item(80, 9).Value = 1997.3334
$ws.Cells.Item(80, 10).Value = 17249.75
$ws.Cells.Item(80, 11).Value = 1997.3334
$ws.Cells.Item(80, 12).Value = 17249.75
$ws.Cells.Item(80, 13).Value = -999.3334
$ws.Cells.Item(80, 14).Value = -19245.75

$ws.Cells.Item(83, 8).Value = 8098.3
$ws.Cells.Item(83, 9).Value = 1997.3334
$ws.Cells.Item(83, 10).Value = 17249.75
$ws.Cells.Item(83, 11).Value = 9986.666999999999
$ws.Cells.Item(83, 12).Value = 86248.75
$ws.Cells.Item(83, 13).Value = -4994.666999999999
$ws.Cells.Item(83, 14).Value = -96232.75

$ws.Cells.Item(97, 8).Value = 676.13635
$ws.Cells.Item(97, 9).Value = 782.1875
$ws.Cells.Item(97, 10).Value = 393.33334
$ws.Cells.Item(97, 11).Value = 782.1875
$ws.Cells.Item(97, 12).Value = 393.33334
$ws.Cells.Item(97, 13).Value = -286.1875
$ws.Cells.Item(97, 14).Value = -1385.33334

$ws.Cells.Item(122, 8).Value = 3546.7368
$ws.Cells.Item(122, 9).Value = 3403.5386
$ws.Cells.Item(122, 11).Value = 10210.6158
$ws.Cells.Item(122, 13).Value = -7760.6158

$ws.Cells.Item(126, 8).Value = 4218.25
$ws.Cells.Item(126, 10).Value = 6749.75
$ws.Cells.Item(126, 12).Value = 20249.25
$ws.Cells.Item(126, 14).Value = -25189.25


$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3035.4167
$ws.Cells.Item(46, 9).Value = 802.5
$ws.Cells.Item(46, 10).Value = 3482
$ws.Cells.Item(46, 11).Value = 802.5
$ws.Cells.Item(46, 12).Value = 3482
$ws.Cells.Item(46, 13).Value = -614.5
$ws.Cells.Item(46, 14).Value = -3858

$ws.Cells.Item(55, 8).Value = 1126
$ws.Cells.Item(55, 10).Value = 2136.1667
$ws.Cells.Item(55, 12).Value = 2136.1667
$ws.Cells.Item(55, 14).Value = -2482.1667

$ws.Cells.Item(136, 8).Value = 4477.381
$ws.Cells.Item(136, 9).Value = 4005.2354
$ws.Cells.Item(136, 10).Value = 6484
$ws.Cells.Item(136, 11).Value = 12015.7062
$ws.Cells.Item(136, 12).Value = 19452
$ws.Cells.Item(136, 13).Value = -9465.706200000001
$ws.Cells.Item(136, 14).Value = -24552


$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1676.5186
$ws.Cells.Item(132, 9).Value = 1213.4
$ws.Cells.Item(132, 11).Value = 3640.2
$ws.Cells.Item(132, 13).Value = -1110.2
